$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ConstantProp")
$rng = $ws.Range("A2:R47")

$vals = @("1","11a","2","3","4","5","8a","8b","9a","--A0 folder","-AlMoO3 folder","AlMoO3-1","AlMoO3-10","AlMoO3-11","AlMoO3-12","AlMoO3-13","AlMoO3-14","AlMoO3-15","AlMoO3-16","AlMoO3-17","AlMoO3-18","AlMoO3-19","AlMoO3-2","AlMoO3-20","AlMoO3-3","AlMoO3-4","AlMoO3-5","AlMoO3-6","AlMoO3-7","AlMoO3-8","AlMoO3-9","Conserv folder","--Mesh folder","")

$rng.AutoFilter(1, $vals, 7)
$rng.AutoFilter(10, @("70"), 7)
$rng.AutoFilter(11, @("4.89E+06"), 7)
Write-Host "done"
